# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Update "OFF" sheet row 2 (H) with new Week 17 totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 559
$wsOff.Range("C2").Value = 388
$wsOff.Range("D2").Value = 116
$wsOff.Range("E2").Value = 59

# Update "DEF" sheet row 2 (H) with new Week 17 totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 478
$wsDef.Range("C2").Value = 329
$wsDef.Range("D2").Value = 114
$wsDef.Range("E2").Value = 49
